# TAC-3831 Fix enable TMS to import trips and edit lists of sub category in excel file
#
# - Adds a new hidden "Lists" sheet (sub-category reference list) after Sheet1
# - Writes "Dry goods" into Sheet1!C2 (first data row default sub category)
# - Wires a list-type data validation on Sheet1!C2:C1048576 sourced from Lists!A3:A30
# - Applies a distinct font color to the "Packed food" helper row on the Lists sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Sheet1: seed the new "Goods Sub Category" default value on row 2
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value = "Dry goods"

# ---------------------------------------------------------------------------
# 2. Add the hidden "Lists" worksheet right after "Sheet1"
# ---------------------------------------------------------------------------
$lists = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$lists.Name = "Lists"

$subCategories = @(
  "Good Sub Category",
  "Packed food ",
  "Diary product ",
  "Beverages ",
  "Fresh food ",
  "Grains & beans ",
  "Animal food ",
  "Household electronics ",
  "General electronics ",
  "Mobiles ",
  "Furniture ",
  "Textiles ",
  "Cosmetics ",
  "Medicine ",
  "Medical equipments ",
  "Medical consumables ",
  "Petrochemicals - Dry ",
  "Petrochemicals - Liquid ",
  "Cars",
  "Tiers ",
  "Spare parts",
  "Lubricants",
  "Steel",
  "Minerals ",
  "Chemicals",
  "Dry goods",
  "Others ",
  "Container",
  "Detergents"
)

for ($i = 0; $i -lt $subCategories.Length; $i++) {
  $row = $i + 1
  $lists.Cells.Item($row, 1).Value = $subCategories[$i]
}

# Row 2 ("Packed food ") carries the distinct dark-grey font used as a group header
$lists.Range("A2").Font.Color = 4473924

$lists.Columns.Item(1).ColumnWidth = 21.140625

$lists.Visible = $false

# ---------------------------------------------------------------------------
# 3. Sheet1 data validation: drop the old blank C1:C1048576 rule (fold C1 into
#    the generic blank rule) and add a list validation sourced from Lists!A3:A30
# ---------------------------------------------------------------------------
$ws1.Range("C1:C1048576").Validation.Delete()
$ws1.Range("C1").Validation.Add(0, 1, 1, [System.Reflection.Missing]::Value)

$ws1.Range("C2:C1048576").Validation.Add(3, 1, 1, "=Lists!A3:A30")

# ---------------------------------------------------------------------------
# 4. Restore the view state captured in the source workbook
# ---------------------------------------------------------------------------
$ws1.Range("I1").Select()
$ws1.Range("R2").Select()

$lists.Range("A10").Select()
$lists.Range("B4").Select()

$ws1.Activate()
